$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data area (rows 2-29, columns B:D) then rebuild only
# the cells that survive in the new layout. Column A (port numbers /
# letters) stays untouched.
$ws.Range("B2:D21").ClearContents()
$ws.Range("B22:D29").ClearContents()

# Row 11 - Distance sensor (unchanged content, but shared-string index
# shifts because of table cleanup)
$ws.Range("B11").Value = "Distance"
$ws.Range("C11").Value = "goal_sense"
$ws.Range("D11").Value = "Distance Sensor"

# Rows 14,16,19,20 - renamed drivetrain motors, NAME (code) column removed
$ws.Range("B14").Value = "Motor"
$ws.Range("D14").Value = "Right Rear Motor"

$ws.Range("B16").Value = "Motor"
$ws.Range("D16").Value = "Right Front Motor"

$ws.Range("B19").Value = "Motor"
$ws.Range("D19").Value = "Left Rear Motor"

$ws.Range("B20").Value = "Motor"
$ws.Range("D20").Value = "Left Front Motor"

# Rows 24-27 - encoders (unchanged content)
$ws.Range("B24").Value = "Encoder (1)"
$ws.Range("C24").Value = "left_enc"
$ws.Range("D24").Value = "Left Custom Encoder"

$ws.Range("B25").Value = "Encoder (2)"
$ws.Range("C25").Value = "left_enc"
$ws.Range("D25").Value = "Left Custom Encoder"

$ws.Range("B26").Value = "Encoder (1)"
$ws.Range("C26").Value = "right_enc"
$ws.Range("D26").Value = "Right Custom Encoder"

$ws.Range("B27").Value = "Encoder (2)"
$ws.Range("C27").Value = "right_enc"
$ws.Range("D27").Value = "Right Custom Encoder"

# Row 28 - new gyro entry (letter port G)
$ws.Range("B28").Value = "Old Gyro"
$ws.Range("D28").Value = "3-Wire Gyro"

# Update the view: zoom level and active selection
$excel.ActiveWindow.Zoom = 160
$ws.Range("D10").Select()
